# S2M12019_Adicionales.xlsx — "Add files via upload"
#
# The sheet (Hoja1) is an attendance/points roster: column B has student
# names (rows 4-30), columns C..BM are daily point entries, and column BN
# sums each row (SUM(C:BM)). This edit fills in a handful of previously
# empty daily-point cells (columns E, K, M, Z, AR) for various students and
# lets the BN totals recompute. Row 13's M cell is entered as a formula
# (=3+3) rather than a literal. Rows 26-30 previously had no BN formula at
# all, so those are (re)created. Finally the active selection on the sheet
# moves from B4 to Z14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---
$ws.Range("E4").Value = 2
$ws.Range("K4").Value = 3

# --- Row 5 ---
$ws.Range("E5").Value = 2
$ws.Range("K5").Value = 3

# --- Row 6 ---
$ws.Range("E6").Value = -1
$ws.Range("K6").Value = 3

# --- Row 7 ---
$ws.Range("E7").Value = 2

# --- Row 8 ---
$ws.Range("E8").Value = 0
$ws.Range("K8").Value = 3

# --- Row 9 ---
$ws.Range("E9").Value = 1
$ws.Range("M9").Value = 3

# --- Row 10 ---
$ws.Range("E10").Value = -1
$ws.Range("M10").Value = 3

# --- Row 11 ---
$ws.Range("E11").Value = -1
$ws.Range("M11").Value = 3

# --- Row 12 ---
$ws.Range("E12").Value = -1
$ws.Range("M12").Value = 3

# --- Row 13 --- (M13 becomes a real formula, not a literal)
$ws.Range("E13").Value = 1
$ws.Range("M13").Formula = "=3+3"
$ws.Range("Z13").Value = 5

# --- Row 14 ---
$ws.Range("E14").Value = -1

# --- Row 15 ---
$ws.Range("E15").Value = -1

# --- Row 16 ---
$ws.Range("E16").Value = 2

# --- Row 18 ---
$ws.Range("E18").Value = -1

# --- Row 19 ---
$ws.Range("E19").Value = 1

# --- Row 20 ---
$ws.Range("E20").Value = -1

# --- Row 21 ---
$ws.Range("E21").Value = 1

# --- Row 22 ---
$ws.Range("E22").Value = -1

# --- Row 23 ---
$ws.Range("E23").Value = 0

# --- Row 24 ---
$ws.Range("E24").Value = 0

# --- Row 25 ---
$ws.Range("E25").Value = -1

# --- Row 26 --- (BN26 previously had no formula at all)
$ws.Range("E26").Value = 1
$ws.Range("BN26").Formula = "=SUM(C26:BM26)"

# --- Row 27 --- (BN27 previously had no formula at all)
$ws.Range("E27").Value = 1
$ws.Range("BN27").Formula = "=SUM(C27:BM27)"

# --- Row 28 --- (BN28 previously had no formula at all)
$ws.Range("E28").Value = -1
$ws.Range("BN28").Formula = "=SUM(C28:BM28)"

# --- Row 29 --- (BN29 previously had no formula at all)
$ws.Range("E29").Value = 2
$ws.Range("AR29").Value = 3
$ws.Range("BN29").Formula = "=SUM(C29:BM29)"

# --- Row 30 --- (BN30 previously had no formula at all)
$ws.Range("E30").Value = -1
$ws.Range("M30").Value = 3
$ws.Range("BN30").Formula = "=SUM(C30:BM30)"

# Move the active selection on the frozen (bottom-right) pane from B4 to Z14.
$ws.Range("Z14").Select()
